$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.714.55"
$ws.Cells.Item(2, 5).Value = "  -2.00%  "
$ws.Cells.Item(3, 4).Value = "3.937.83"
$ws.Cells.Item(3, 5).Value = "  -2.38%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).Value = "'535.69"
$ws.Cells.Item(5, 5).Value = "  +2.78%  "
$ws.Cells.Item(6, 4).Value = "'147.82"
$ws.Cells.Item(6, 5).Value = "  +0.45%  "
$ws.Cells.Item(7, 4).Value = "3.935.12"
$ws.Cells.Item(7, 5).Value = "  -2.20%  "
$ws.Cells.Item(8, 5).Value = "  -4.51%  "
$ws.Cells.Item(9, 4).Value = "'1.00"
$ws.Cells.Item(9, 5).Value = "  +0.02%  "
$ws.Cells.Item(10, 5).Value = "  -4.68%  "
$ws.Cells.Item(11, 5).Value = "  -5.81%  "
$ws.Cells.Item(12, 4).Value = "'55.30"
$ws.Cells.Item(12, 5).Value = "  +12.24%  "
$ws.Cells.Item(13, 4).Value = "'0.0000314"
$ws.Cells.Item(13, 5).Value = "  -3.96%  "
$ws.Cells.Item(14, 4).Value = "'10.56"
$ws.Cells.Item(14, 5).Value = "  -5.22%  "
$ws.Cells.Item(15, 4).Value = "4.578.90"
$ws.Cells.Item(15, 5).Value = "  -2.12%  "
$ws.Cells.Item(16, 4).Value = "3.949.69"
$ws.Cells.Item(16, 5).Value = "  -2.30%  "
$ws.Cells.Item(17, 4).Value = "'20.54"
$ws.Cells.Item(17, 5).Value = "  -3.39%  "
$ws.Cells.Item(18, 4).Value = "'13.79"
$ws.Cells.Item(18, 5).Value = "  -3.21%  "
$ws.Cells.Item(19, 5).Value = "  -1.61%  "
$ws.Cells.Item(20, 5).Value = "  -4.60%  "
$ws.Cells.Item(21, 4).Value = "70.837.75"
$ws.Cells.Item(21, 5).Value = "  -1.73%  "
$ws.Cells.Item(22, 4).Value = "'422.94"
$ws.Cells.Item(22, 5).Value = "  -4.78%  "
$ws.Cells.Item(23, 4).Value = "'3.57"
$ws.Cells.Item(23, 5).Value = "  -0.69%  "
$ws.Cells.Item(24, 4).Value = "'96.68"
$ws.Cells.Item(24, 5).Value = "  -8.21%  "
$ws.Cells.Item(25, 5).Value = "  +4.24%  "
$ws.Cells.Item(26, 4).Value = "'14.33"
$ws.Cells.Item(26, 5).Value = "  -5.71%  "
$ws.Cells.Item(27, 4).Value = "'11.32"
$ws.Cells.Item(27, 5).Value = "  -2.18%  "
$ws.Cells.Item(28, 4).Value = "'3.78"
$ws.Cells.Item(28, 5).Value = "  +15.51%  "
$ws.Cells.Item(29, 4).Value = "'10.59"
$ws.Cells.Item(29, 5).Value = "  -4.36%  "
$ws.Cells.Item(30, 4).Value = "'5.86"
$ws.Cells.Item(30, 5).Value = "  +0.76%  "
$ws.Cells.Item(31, 4).Value = "'36.25"
$ws.Cells.Item(31, 5).Value = "  -4.15%  "
$ws.Cells.Item(32, 4).Value = "'7.76"
$ws.Cells.Item(32, 5).Value = "  +14.52%  "
$ws.Cells.Item(33, 4).Value = "'50.77"
$ws.Cells.Item(33, 5).Value = "  +19.46%  "
$ws.Cells.Item(34, 2).Value = "Hedera"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(34, 4).Value = "'0.130"
$ws.Cells.Item(34, 5).Value = "  -0.93%  "
$ws.Cells.Item(35, 2).Value = "Cosmos"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(35, 4).Value = "'13.29"
$ws.Cells.Item(35, 5).Value = "  -3.44%  "
$ws.Cells.Item(36, 4).Value = "'684.05"
$ws.Cells.Item(36, 5).Value = "  +1.21%  "
$ws.Cells.Item(37, 4).Value = "'65.01"
$ws.Cells.Item(37, 5).Value = "  -3.95%  "
$ws.Cells.Item(38, 4).Value = "'0.437"
$ws.Cells.Item(38, 5).Value = "  +2.27%  "
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).Value = "'0.149"
$ws.Cells.Item(39, 5).Value = "  -1.69%  "
$ws.Cells.Item(40, 2).Value = "PEPE"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(40, 4).Value = "0.0₃0815"
$ws.Cells.Item(40, 5).Value = "  -5.31%  "
$ws.Cells.Item(41, 4).Value = "'3.36"
$ws.Cells.Item(41, 5).Value = "  -4.41%  "
$ws.Cells.Item(42, 4).Value = "'0.999"
$ws.Cells.Item(42, 5).Value = "  +0.06%  "
$ws.Cells.Item(43, 5).Value = "  +0.15%  "
$ws.Cells.Item(44, 5).Value = "  -4.32%  "
$ws.Cells.Item(45, 5).Value = "  -0.96%  "
$ws.Cells.Item(46, 2).Value = "THORChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(46, 4).Value = "'9.93"
$ws.Cells.Item(46, 5).Value = "  +4.88%  "
$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(47, 4).Value = "'0.148"
$ws.Cells.Item(47, 5).Value = "  -6.43%  "
$ws.Cells.Item(48, 4).Value = "'2.68"
$ws.Cells.Item(48, 5).Value = "  -1.34%  "
$ws.Cells.Item(49, 4).Value = "'3.35"
$ws.Cells.Item(49, 5).Value = "  -3.95%  "
$ws.Cells.Item(50, 4).Value = "'2.99"
$ws.Cells.Item(50, 5).Value = "  -2.66%  "
$ws.Cells.Item(51, 2).Value = "Monero"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51, 4).Value = "'144.70"
$ws.Cells.Item(51, 5).Value = "  -0.52%  "
